# cap nhat bo sung tinh luong
# Insert two new salary-bracket columns (70, 80) before the old "Ti le"/"Bat cap"
# columns, fill in the new percentage progression, refresh the "*" formula
# strings for a few employees, and tidy up the leftover wrap-text style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two columns at J:K - this shifts the existing J (Ti le) and K (Bat
#    cap) columns to L and M, carrying their values/styles/formats along.
$ws.Range("J1:K1").EntireColumn.Insert()

# 2) New header values for the inserted bracket columns.
$ws.Cells.Item(1, 10).Value = 70
$ws.Cells.Item(1, 11).Value = 80

# 3) Row 2 (Minh) - fill the new bracket columns and update the cap ratio.
$ws.Cells.Item(2, 10).Value = 0.01
$ws.Cells.Item(2, 11).Value = 0.01
$ws.Cells.Item(2, 12).Value = 0.3

# 4) Rows 3 & 4 (Hai, Dung) share the same updated percentage progression and
#    gain two more brackets.
foreach ($r in 3, 4) {
    $ws.Cells.Item($r, 5).Value = 0.013
    $ws.Cells.Item($r, 6).Value = 0.015
    $ws.Cells.Item($r, 7).Value = 0.016
    $ws.Cells.Item($r, 8).ClearFormats()
    $ws.Cells.Item($r, 8).Value = 0.017
    $ws.Cells.Item($r, 9).Value = 0.018
    $ws.Cells.Item($r, 10).Value = 0.019
    $ws.Cells.Item($r, 11).Value = 0.02
}

# Dung (row 4) also gets the combined "*|Minh:0.7|Duong:0.7" cap formula.
$ws.Cells.Item(4, 13).Value = "*|Minh:0.7|Duong:0.7"

# 5) Rows 9-12 (Tuan, Nam, Dong, Hong) share another updated progression.
foreach ($r in 9, 10, 11, 12) {
    $ws.Cells.Item($r, 5).Value = 0.012
    $ws.Cells.Item($r, 6).Value = 0.014
    $ws.Cells.Item($r, 7).Value = 0.015
    $ws.Cells.Item($r, 8).Value = 0.016
    $ws.Cells.Item($r, 9).Value = 0.017
    $ws.Cells.Item($r, 10).Value = 0.018
    $ws.Cells.Item($r, 11).Value = 0.019
}

# Tuan, Nam and Hong now reference the combined cap formula too (Dong keeps
# the plain "*").
$ws.Cells.Item(9, 13).Value = "*|Minh:0.7|Duong:0.7"
$ws.Cells.Item(10, 13).Value = "*|Minh:0.7|Duong:0.7"
$ws.Cells.Item(12, 13).Value = "*|Minh:0.7|Duong:0.7"

# 6) Column widths for the newly inserted bracket columns and the wider cap
#    column.
$ws.Range("J1:K1").ColumnWidth = 5.166666666666667
$ws.Range("M1").ColumnWidth = 20.022135416666668

# 7) Restore the selection to the last edited cell.
$ws.Range("M11").Select()
